$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate the K column (G) to use a constant value of 1 instead of the
# previous "Strike#" derived values for data rows 3-7.
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 1
